# Apply the Dec 14 2023 05:56:04 UTC cryptos-list refresh (prices / 1h-volume
# deltas, plus a few rank swaps) to Sheet1, cell by cell, exactly as the diff
# describes it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell. Numeric-looking strings (e.g. "249.42")
# would otherwise be auto-coerced to a Number by Excel on assignment, which
# would lose the source's exact text formatting (trailing zeros, "249.44" vs
# 249.44, double-dot big numbers, etc). A leading apostrophe is the normal
# Excel idiom for "store this as text" and is stripped from the stored text,
# leaving only the cell marked as text (quote-prefixed).
function Set-TextValue($range, $value) {
    $looksNumeric = $false
    try {
        [void][double]::Parse($value.Trim(), [System.Globalization.CultureInfo]::InvariantCulture)
        $looksNumeric = $true
    } catch {
        $looksNumeric = $false
    }
    if ($looksNumeric) {
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

# Row 2: Bitcoin
Set-TextValue $ws.Range('D2') '42.792.35'
Set-TextValue $ws.Range('E2') '  +4.33%  '

# Row 3: Ethereum
Set-TextValue $ws.Range('D3') '2.259.09'
Set-TextValue $ws.Range('E3') '  +3.99%  '

# Row 4: TetherUSD
Set-TextValue $ws.Range('E4') '  -0.01%  '

# Row 5: BNB
Set-TextValue $ws.Range('D5') '249.42'
Set-TextValue $ws.Range('E5') '  +1.15%  '

# Row 6: XRP
Set-TextValue $ws.Range('D6') '0.625'
Set-TextValue $ws.Range('E6') '  +1.38%  '

# Row 7: Solana
Set-TextValue $ws.Range('D7') '71.41'
Set-TextValue $ws.Range('E7') '  +7.91%  '

# Row 8: USDC
Set-TextValue $ws.Range('E8') '  -0.09%  '

# Row 9: Cardano
Set-TextValue $ws.Range('D9') '0.650'
Set-TextValue $ws.Range('E9') '  +15.06%  '

# Row 10: Avalanche
Set-TextValue $ws.Range('D10') '38.36'
Set-TextValue $ws.Range('E10') '  +7.94%  '

# Row 11: Dogecoin
Set-TextValue $ws.Range('B11') 'Dogecoin'
Set-TextValue $ws.Range('C11') 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range('D11') '0.0976'
Set-TextValue $ws.Range('E11') '  +5.32%  '

# Row 12: OKB
Set-TextValue $ws.Range('B12') 'OKB'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D12') '59.45'
Set-TextValue $ws.Range('E12') '  -0.62%  '

# Row 13: Polkadot
Set-TextValue $ws.Range('D13') '7.37'
Set-TextValue $ws.Range('E13') '  +7.64%  '

# Row 14: TRON
Set-TextValue $ws.Range('E14') '  +1.64%  '

# Row 15: WrappedliquidstakedEther2.0
Set-TextValue $ws.Range('D15') '2.595.56'
Set-TextValue $ws.Range('E15') '  +4.16%  '

# Row 16: Chainlink
Set-TextValue $ws.Range('D16') '14.87'
Set-TextValue $ws.Range('E16') '  +4.61%  '

# Row 17: Polygon
Set-TextValue $ws.Range('D17') '0.881'
Set-TextValue $ws.Range('E17') '  +3.11%  '

# Row 18: WrappedEther
Set-TextValue $ws.Range('D18') '2.273.11'
Set-TextValue $ws.Range('E18') '  +4.94%  '

# Row 19: WrappedBTC
Set-TextValue $ws.Range('D19') '42.734.12'
Set-TextValue $ws.Range('E19') '  +4.35%  '

# Row 20: ShibaInu
Set-TextValue $ws.Range('D20') '0.0000101'
Set-TextValue $ws.Range('E20') '  +7.89%  '

# Row 21: Uniswap
Set-TextValue $ws.Range('D21') '6.32'
Set-TextValue $ws.Range('E21') '  +3.88%  '

# Row 22: Litecoin
Set-TextValue $ws.Range('D22') '73.08'
Set-TextValue $ws.Range('E22') '  +2.41%  '

# Row 23: BitcoinCash
Set-TextValue $ws.Range('D23') '235.65'
Set-TextValue $ws.Range('E23') '  +2.78%  '

# Row 24: ImmutableX
Set-TextValue $ws.Range('D24') '2.10'
Set-TextValue $ws.Range('E24') '  +0.68%  '

# Row 25: WEMIXToken
Set-TextValue $ws.Range('D25') '3.95'
Set-TextValue $ws.Range('E25') '  +7.58%  '

# Row 26: Cosmos
Set-TextValue $ws.Range('D26') '11.47'
Set-TextValue $ws.Range('E26') '  +1.24%  '

# Row 27: Dai
Set-TextValue $ws.Range('E27') '  -0.08%  '

# Row 28: PancakeSwap
Set-TextValue $ws.Range('E28') '  +0.80%  '

# Row 29: LEO
Set-TextValue $ws.Range('E29') '  -1.61%  '

# Row 30: EthereumClassic
Set-TextValue $ws.Range('B30') 'EthereumClassic'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D30') '21.55'
Set-TextValue $ws.Range('E30') '  +6.83%  '

# Row 31: Toncoin
Set-TextValue $ws.Range('B31') 'Toncoin'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D31') '2.10'
Set-TextValue $ws.Range('E31') '  +5.04%  '

# Row 32: Monero
Set-TextValue $ws.Range('B32') 'Monero'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D32') '167.65'
Set-TextValue $ws.Range('E32') '  -0.74%  '

# Row 33: InternetComputer(DFINITY)
Set-TextValue $ws.Range('D33') '6.53'
Set-TextValue $ws.Range('E33') '  +16.22%  '

# Row 34: Kaspa
Set-TextValue $ws.Range('E34') '  +5.61%  '

# Row 35: Hedera
Set-TextValue $ws.Range('B35') 'Hedera'
Set-TextValue $ws.Range('C35') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D35') '0.0800'
Set-TextValue $ws.Range('E35') '  +6.56%  '

# Row 36: InjectiveProtocol
Set-TextValue $ws.Range('B36') 'InjectiveProtocol'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D36') '31.45'
Set-TextValue $ws.Range('E36') '  +28.07%  '

# Row 37: Stellar
Set-TextValue $ws.Range('E37') '  +3.83%  '

# Row 38: RenderToken
Set-TextValue $ws.Range('D38') '4.46'
Set-TextValue $ws.Range('E38') '  +12.47%  '

# Row 39: Filecoin
Set-TextValue $ws.Range('D39') '4.74'
Set-TextValue $ws.Range('E39') '  +4.17%  '

# Row 40: VeChain
Set-TextValue $ws.Range('D40') '0.0322'
Set-TextValue $ws.Range('E40') '  +5.90%  '

# Row 41: LidoDAOToken
Set-TextValue $ws.Range('D41') '2.32'
Set-TextValue $ws.Range('E41') '  +6.07%  '

# Row 42: Celestia
Set-TextValue $ws.Range('D42') '12.80'
Set-TextValue $ws.Range('E42') '  +12.92%  '

# Row 43: THORChain
Set-TextValue $ws.Range('D43') '5.81'
Set-TextValue $ws.Range('E43') '  +6.27%  '

# Row 44: FraxShare
Set-TextValue $ws.Range('D44') '9.20'
Set-TextValue $ws.Range('E44') '  +9.37%  '

# Row 45: MultiversX
Set-TextValue $ws.Range('D45') '62.12'
Set-TextValue $ws.Range('E45') '  +3.05%  '

# Row 46: Algorand
Set-TextValue $ws.Range('E46') '  +5.46%  '

# Row 47: FTXToken
Set-TextValue $ws.Range('D47') '4.84'
Set-TextValue $ws.Range('E47') '  -0.68%  '

# Row 48: Cronos
Set-TextValue $ws.Range('E48') '  +3.44%  '

# Row 49: BinanceUSD
Set-TextValue $ws.Range('E49') '  +0.05%  '

# Row 50: ARBITRUM
Set-TextValue $ws.Range('D50') '1.16'
Set-TextValue $ws.Range('E50') '  +1.44%  '

# Row 51: TrustWalletToken
Set-TextValue $ws.Range('E51') '  +3.96%  '

